$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5716.8945
$ws.Range("M33").Value = -6278.75
$ws.Range("I33").Value = 6507.75
$ws.Range("K33").Value = 6507.75
$ws.Range("H49").Value = 641.5
$ws.Range("I49").Value = 641.5
$ws.Range("K49").Value = 1924.5
$ws.Range("M49").Value = -1788.5
$ws.Range("I61").Value = 126
$ws.Range("K61").Value = 378
$ws.Range("H61").Value = 1344.25
$ws.Range("M61").Value = -206
$ws.Range("I62").Value = 24459.75
$ws.Range("L62").Value = 9524.571
$ws.Range("N62").Value = -10772.571
$ws.Range("H62").Value = 12843.5
$ws.Range("K62").Value = 24459.75
$ws.Range("M62").Value = -23835.75
$ws.Range("J62").Value = 9524.571
$ws.Range("L65").Value = 47622.855
$ws.Range("H65").Value = 12843.5
$ws.Range("I65").Value = 24459.75
$ws.Range("J65").Value = 9524.571
$ws.Range("M65").Value = -119178.75
$ws.Range("N65").Value = -53862.855
$ws.Range("K65").Value = 122298.75
$ws.Range("L69").Value = 40392.858
$ws.Range("K69").Value = 0
$ws.Range("N69").Value = -42140.858
$ws.Range("M69").ClearContents()
$ws.Range("J69").Value = 13464.286
$ws.Range("H69").Value = 13464.286
$ws.Range("I69").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -129914.574
$ws.Range("K72").Value = 0
$ws.Range("H72").Value = 13464.286
$ws.Range("L72").Value = 121178.574
$ws.Range("J72").Value = 13464.286
$ws.Range("I72").Value = 0
$ws.Range("M82").Value = -11049.9095
$ws.Range("H82").Value = 4961.846
$ws.Range("I82").Value = 3818.6365
$ws.Range("K82").Value = 11455.9095
$ws.Range("M85").Value = -10051.9095
$ws.Range("K85").Value = 11455.9095
$ws.Range("H85").Value = 4961.846
$ws.Range("I85").Value = 3818.6365
$ws.Range("L96").Value = 7798.5
$ws.Range("K96").Value = 598.5
$ws.Range("M96").Value = 774.5
$ws.Range("J96").Value = 2599.5
$ws.Range("H96").Value = 1999.5
$ws.Range("N96").Value = -10544.5
$ws.Range("I96").Value = 199.5
$ws.Range("H98").Value = 1369.1177
$ws.Range("K98").Value = 1153.0769
$ws.Range("I98").Value = 1153.0769
$ws.Range("M98").Value = 344.9231
$ws.Range("I103").Value = 1625.25
$ws.Range("H103").Value = 1695.909
$ws.Range("K103").Value = 4875.75
$ws.Range("M103").Value = -4289.75
$ws.Range("H106").Value = 9150
$ws.Range("L106").Value = 15691.667
$ws.Range("K106").Value = 7187.5
$ws.Range("N106").Value = -16953.667
$ws.Range("J106").Value = 15691.667
$ws.Range("M106").Value = -6556.5
$ws.Range("I106").Value = 7187.5
$ws.Range("H107").Value = 1976
$ws.Range("K107").Value = 1408
$ws.Range("M107").Value = 512
$ws.Range("I107").Value = 1408
$ws.Range("M113").Value = 146.1109999999999
$ws.Range("H113").Value = 4921.731
$ws.Range("K113").Value = 3107.889
$ws.Range("I113").Value = 3107.889
$ws.Range("I122").Value = 1153.0769
$ws.Range("K122").Value = 3459.2307
$ws.Range("H122").Value = 1369.1177
$ws.Range("M122").Value = -1009.2307
$ws.Range("I131").Value = 1218.0555
$ws.Range("M131").Value = 1385.8335
$ws.Range("K131").Value = 3654.1665
$ws.Range("H131").Value = 2206.5789
$ws.Range("K132").Value = 22866.429
$ws.Range("H132").Value = 8749.276
$ws.Range("M132").Value = -20336.429
$ws.Range("I132").Value = 7622.143
$ws.Range("M137").Value = -67.875
$ws.Range("N137").Value = -15743.25
$ws.Range("J137").Value = 3547.75
$ws.Range("I137").Value = 872.625
$ws.Range("H137").Value = 2210.1875
$ws.Range("L137").Value = 10643.25
$ws.Range("K137").Value = 2617.875
$ws.Range("H138").Value = 4419.476
$ws.Range("K138").Value = 7980.1428
$ws.Range("L138").Value = 18536.715
$ws.Range("M138").Value = -2840.1428
$ws.Range("J138").Value = 6178.905
$ws.Range("I138").Value = 2660.0476
$ws.Range("N138").Value = -28816.715

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L45").Value = 3996.6155
$ws.Range("H45").Value = 2892.1875
$ws.Range("N45").Value = -4750.6155
$ws.Range("J45").Value = 3996.6155
$ws.Range("K53").Value = 12399.4
$ws.Range("I53").Value = 12399.4
$ws.Range("H53").Value = 16999.334
$ws.Range("M53").Value = -11717.4
$ws.Range("L62").Value = 180000
$ws.Range("N62").Value = -181248
$ws.Range("H62").Value = 180000
$ws.Range("J62").Value = 180000
$ws.Range("K63").Value = 2850
$ws.Range("H63").Value = 7141.6665
$ws.Range("M63").Value = -2164
$ws.Range("I63").Value = 2850
$ws.Range("N64").Value = -150496
$ws.Range("J64").Value = 150000
$ws.Range("H64").Value = 150000
$ws.Range("L64").Value = 150000
$ws.Range("L65").Value = 540000
$ws.Range("H65").Value = 180000
$ws.Range("J65").Value = 180000
$ws.Range("N65").Value = -546240
$ws.Range("M66").Value = -10818
$ws.Range("I66").Value = 2850
$ws.Range("K66").Value = 14250
$ws.Range("H66").Value = 7141.6665
$ws.Range("J67").Value = 150000
$ws.Range("H67").Value = 150000
$ws.Range("L67").Value = 150000
$ws.Range("N67").Value = -151716
$ws.Range("I74").Value = 26924.316
$ws.Range("K74").Value = 26924.316
$ws.Range("M74").Value = -26050.316
$ws.Range("H74").Value = 26924.316
$ws.Range("H76").Value = 55100
$ws.Range("J76").Value = 55100
$ws.Range("N76").Value = -55776
$ws.Range("L76").Value = 55100
$ws.Range("I77").Value = 26924.316
$ws.Range("H77").Value = 26924.316
$ws.Range("M77").Value = -130253.58
$ws.Range("K77").Value = 134621.58
$ws.Range("J79").Value = 55100
$ws.Range("H79").Value = 55100
$ws.Range("L79").Value = 55100
$ws.Range("N79").Value = -57440
$ws.Range("N97").Value = -6966.375
$ws.Range("J97").Value = 5974.375
$ws.Range("L97").Value = 5974.375
$ws.Range("K97").Value = 5778.35
$ws.Range("I97").Value = 5778.35
$ws.Range("M97").Value = -5282.35
$ws.Range("H97").Value = 5834.357
$ws.Range("H110").Value = 3533.6365
$ws.Range("I110").Value = 3207.889
$ws.Range("N110").Value = -9089.5
$ws.Range("L110").Value = 4999.5
$ws.Range("K110").Value = 3207.889
$ws.Range("J110").Value = 4999.5
$ws.Range("M110").Value = -1162.889

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M20").Value = -3501.8333
$ws.Range("J20").Value = 3963
$ws.Range("I20").Value = 3748.8333
$ws.Range("H20").Value = 3871.2144
$ws.Range("K20").Value = 3748.8333
$ws.Range("L20").Value = 3963
$ws.Range("N20").Value = -4457
$ws.Range("L94").Value = 4246.5
$ws.Range("I94").Value = 2485.1428
$ws.Range("H94").Value = 3125.6365
$ws.Range("K94").Value = 2485.1428
$ws.Range("M94").Value = -2034.1428
$ws.Range("N94").Value = -5148.5
$ws.Range("J94").Value = 4246.5
$ws.Range("M134").Value = -3836.099999999999
$ws.Range("I134").Value = 2123.7
$ws.Range("H134").Value = 2439.394
$ws.Range("K134").Value = 6371.099999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2421.6572
$ws.Range("K31").Value = 2148.7036
$ws.Range("N31").Value = -3932.875
$ws.Range("M31").Value = -1853.7036
$ws.Range("J31").Value = 3342.875
$ws.Range("L31").Value = 3342.875
$ws.Range("I31").Value = 2148.7036
$ws.Range("L34").Value = 3342.875
$ws.Range("N34").Value = -3746.875
$ws.Range("M34").Value = -1946.7036
$ws.Range("K34").Value = 2148.7036
$ws.Range("J34").Value = 3342.875
$ws.Range("H34").Value = 2421.6572
$ws.Range("I34").Value = 2148.7036
$ws.Range("J86").Value = 4632
$ws.Range("H86").Value = 53960.125
$ws.Range("M86").Value = -82434
$ws.Range("L86").Value = 4632
$ws.Range("I86").Value = 83557
$ws.Range("N86").Value = -6878
$ws.Range("K86").Value = 83557
$ws.Range("M89").Value = -412169
$ws.Range("L89").Value = 23160
$ws.Range("I89").Value = 83557
$ws.Range("K89").Value = 417785
$ws.Range("N89").Value = -34392
$ws.Range("J89").Value = 4632
$ws.Range("H89").Value = 53960.125
$ws.Range("L96").Value = 17800
$ws.Range("J96").Value = 17800
$ws.Range("H96").Value = 17800
$ws.Range("N96").Value = -23292
$ws.Range("H107").Value = 471.94736
$ws.Range("K107").Value = 453.94116
$ws.Range("M107").Value = 1466.05884
$ws.Range("I107").Value = 453.94116
$ws.Range("J141").Value = 900000
$ws.Range("H141").Value = 900000
$ws.Range("L141").Value = 900000
$ws.Range("N141").Value = -910360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I7").Value = 301
$ws.Range("J7").Value = 5555
$ws.Range("K7").Value = 903
$ws.Range("L7").Value = 16665
$ws.Range("N7").Value = -16889
$ws.Range("H7").Value = 3803.6667
$ws.Range("M7").Value = -791
$ws.Range("K18").Value = 1603.2
$ws.Range("H18").Value = 534.4
$ws.Range("M18").Value = -1434.2
$ws.Range("I18").Value = 534.4
$ws.Range("J39").Value = 9995
$ws.Range("H39").Value = 9995
$ws.Range("L39").Value = 29985
$ws.Range("N39").Value = -30573
$ws.Range("N97").Value = -2591.9999
$ws.Range("J97").Value = 533.3333
$ws.Range("L97").Value = 1599.9999
$ws.Range("H97").Value = 394.5
$ws.Range("H98").Value = 2469.375
$ws.Range("K98").Value = 1197
$ws.Range("I98").Value = 399
$ws.Range("M98").Value = 301
$ws.Range("H104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M114").Value = 137
$ws.Range("I114").Value = 1039
$ws.Range("K114").Value = 3117
$ws.Range("H114").Value = 1423.4
$ws.Range("L121").Value = 2613
$ws.Range("N121").Value = -5233
$ws.Range("H121").Value = 1236.625
$ws.Range("K121").Value = 4563
$ws.Range("J121").Value = 871
$ws.Range("I121").Value = 1521
$ws.Range("M121").Value = -3253
$ws.Range("I131").Value = 22485
$ws.Range("M131").Value = -62415
$ws.Range("K131").Value = 67455
$ws.Range("H131").Value = 5574119.5
$ws.Range("I136").Value = 2455
$ws.Range("M136").Value = -2265
$ws.Range("K136").Value = 7365
$ws.Range("H136").Value = 4220.25
$ws.Range("K138").Value = 1430.1429
$ws.Range("M138").Value = 3709.8571
$ws.Range("I138").Value = 476.7143
$ws.Range("H138").Value = 476.7143
$ws.Range("H140").Value = 2370.9473
$ws.Range("I140").Value = 1929.3529
$ws.Range("K140").Value = 5788.0587
$ws.Range("M140").Value = -608.0587000000005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H52").Value = 39990
$ws.Range("I52").Value = 0
$ws.Range("H53").Value = 49999
$ws.Range("L53").Value = 49999
$ws.Range("J53").Value = 49999
$ws.Range("N53").Value = -51261
$ws.Range("I57").Value = 5055
$ws.Range("K57").Value = 5055
$ws.Range("M57").Value = -4235
$ws.Range("H57").Value = 18821.6
$ws.Range("I70").Value = 13550.2
$ws.Range("M70").Value = -13280.2
$ws.Range("K70").Value = 13550.2
$ws.Range("H70").Value = 13678.571
$ws.Range("I73").Value = 13550.2
$ws.Range("M73").Value = -12614.2
$ws.Range("K73").Value = 13550.2
$ws.Range("H73").Value = 13678.571
$ws.Range("L80").Value = 3499.9285
$ws.Range("M80").Value = -957
$ws.Range("I80").Value = 1955
$ws.Range("H80").Value = 2786.8845
$ws.Range("J80").Value = 3499.9285
$ws.Range("N80").Value = -5495.9285
$ws.Range("K80").Value = 1955
$ws.Range("K83").Value = 9775
$ws.Range("H83").Value = 2786.8845
$ws.Range("I83").Value = 1955
$ws.Range("L83").Value = 17499.6425
$ws.Range("M83").Value = -4783
$ws.Range("N83").Value = -27483.6425
$ws.Range("J83").Value = 3499.9285
$ws.Range("M102").Value = 896.7273
$ws.Range("I102").Value = 725.2727
$ws.Range("J102").Value = 3027
$ws.Range("L102").Value = 3027
$ws.Range("K102").Value = 725.2727
$ws.Range("H102").Value = 1218.5
$ws.Range("N102").Value = -6271
$ws.Range("H107").Value = 44869.78
$ws.Range("K107").Value = 77935.234
$ws.Range("M107").Value = -76015.234
$ws.Range("I107").Value = 77935.234
$ws.Range("M113").Value = -275330
$ws.Range("H113").Value = 352500
$ws.Range("K113").Value = 277500
$ws.Range("I113").Value = 277500
$ws.Range("N122").Value = -14764.6
$ws.Range("I122").Value = 2291.6
$ws.Range("K122").Value = 6874.799999999999
$ws.Range("H122").Value = 2623.8
$ws.Range("J122").Value = 3288.2
$ws.Range("M122").Value = -4424.799999999999
$ws.Range("L122").Value = 9864.599999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M22").Value = -1117.5
$ws.Range("H22").Value = 2672.9333
$ws.Range("K22").Value = 1412.5
$ws.Range("L22").Value = 4113.4287
$ws.Range("I22").Value = 1412.5
$ws.Range("J22").Value = 4113.4287
$ws.Range("N22").Value = -4703.4287
$ws.Range("J27").Value = 4113.4287
$ws.Range("M27").Value = -1305.5
$ws.Range("I27").Value = 1412.5
$ws.Range("H27").Value = 2672.9333
$ws.Range("K27").Value = 1412.5
$ws.Range("N27").Value = -4327.4287
$ws.Range("L27").Value = 4113.4287
$ws.Range("K55").Value = 755.8461
$ws.Range("J55").Value = 1500
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = -582.8461
$ws.Range("H55").Value = 855.06665
$ws.Range("N55").Value = -1846
$ws.Range("I55").Value = 755.8461
$ws.Range("K63").Value = 0
$ws.Range("H63").Value = 63999
$ws.Range("M63").ClearContents()
$ws.Range("I63").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H66").Value = 63999
$ws.Range("K68").Value = 4159.8
$ws.Range("M68").Value = -3410.8
$ws.Range("N68").Value = -6123
$ws.Range("I68").Value = 4159.8
$ws.Range("H68").Value = 4366.5557
$ws.Range("J68").Value = 4625
$ws.Range("L68").Value = 4625
$ws.Range("I71").Value = 4159.8
$ws.Range("K71").Value = 20799
$ws.Range("M71").Value = -17055
$ws.Range("J71").Value = 4625
$ws.Range("N71").Value = -30613
$ws.Range("L71").Value = 23125
$ws.Range("H71").Value = 4366.5557
$ws.Range("I122").Value = 3418.0908
$ws.Range("K122").Value = 10254.2724
$ws.Range("H122").Value = 4137.125
$ws.Range("M122").Value = -7804.2724
$ws.Range("I136").Value = 2156.9524
$ws.Range("M136").Value = -3920.8572
$ws.Range("K136").Value = 6470.8572
$ws.Range("H136").Value = 3132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 5499
$ws.Range("H62").Value = 167639.72
$ws.Range("K62").Value = 5499
$ws.Range("M62").Value = -4875
$ws.Range("H65").Value = 167639.72
$ws.Range("I65").Value = 5499
$ws.Range("M65").Value = -24375
$ws.Range("K65").Value = 27495
$ws.Range("L96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("H96").Value = 2977.111
$ws.Range("N96").Value = -6746
$ws.Range("M113").Value = -1277.7896
$ws.Range("H113").Value = 1158.3043
$ws.Range("K113").Value = 3447.7896
$ws.Range("N113").Value = -7943.75
$ws.Range("L113").Value = 3603.75
$ws.Range("I113").Value = 1149.2632
$ws.Range("J113").Value = 1201.25
$ws.Range("I122").Value = 2137.9644
$ws.Range("K122").Value = 6413.8932
$ws.Range("H122").Value = 2128.1516
$ws.Range("M122").Value = -3963.8932
$ws.Range("K132").Value = 176015.835
$ws.Range("H132").Value = 58671.945
$ws.Range("M132").Value = -173485.835
$ws.Range("I132").Value = 58671.945
